$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8518
$ws1.Range("F4").Value = 382
$ws1.Range("F5").Value = 25

# Update "全部类型" (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8518
$ws4.Range("F4").Value = 382
$ws4.Range("F5").Value = 25
